# Auto-generated Excel COM-interop script
# Applies scheduled-runner price/profit updates across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 201
$ws.Range("I9").Value = 201
$ws.Range("K9").Value = 201
$ws.Range("M9").Value = -32

$ws.Range("H16").Value = 100
$ws.Range("J16").Value = 100
$ws.Range("L16").Value = 100
$ws.Range("N16").Value = -560

$ws.Range("H62").Value = 6298.5
$ws.Range("I62").Value = 3916.4
$ws.Range("K62").Value = 3916.4
$ws.Range("M62").Value = -3292.4

$ws.Range("H65").Value = 6298.5
$ws.Range("I65").Value = 3916.4
$ws.Range("K65").Value = 19582
$ws.Range("M65").Value = -16462

$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

$ws.Range("H76").Value = 5910.8096
$ws.Range("I76").Value = 4648.364
$ws.Range("K76").Value = 4648.364
$ws.Range("M76").Value = -4333.364

$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

$ws.Range("H79").Value = 5910.8096
$ws.Range("I79").Value = 4648.364
$ws.Range("K79").Value = 4648.364
$ws.Range("M79").Value = -3556.364

$ws.Range("H106").Value = 25862.111
$ws.Range("I106").Value = 28682.75
$ws.Range("K106").Value = 28682.75
$ws.Range("M106").Value = -28051.75

$ws.Range("H116").Value = 7870.222
$ws.Range("J116").Value = 8539.799999999999
$ws.Range("L116").Value = 8539.799999999999
$ws.Range("N116").Value = -15423.8

$ws.Range("H132").Value = 1221.7354
$ws.Range("I132").Value = 1184.0312
$ws.Range("K132").Value = 3552.0936
$ws.Range("M132").Value = -1022.0936

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2000
$ws.Range("I61").Value = 2000
$ws.Range("K61").Value = 2000
$ws.Range("M61").Value = -1788

$ws.Range("H63").Value = 4629.08
$ws.Range("I63").Value = 3545.5625
$ws.Range("K63").Value = 3545.5625
$ws.Range("M63").Value = -2859.5625

$ws.Range("H66").Value = 4629.08
$ws.Range("I66").Value = 3545.5625
$ws.Range("K66").Value = 17727.8125
$ws.Range("M66").Value = -14295.8125

$ws.Range("H136").Value = 2000
$ws.Range("I136").Value = 2000
$ws.Range("K136").Value = 6000
$ws.Range("M136").Value = -3450

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1600
$ws.Range("I20").Value = 1600
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 1600
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -1353
$ws.Range("N20").ClearContents()

$ws.Range("H130").Value = 35709
$ws.Range("I130").Value = 35709
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 35709
$ws.Range("L130").Value = 0
$ws.Range("M130").Value = -30689
$ws.Range("N130").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1166.6666
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()

$ws.Range("H141").Value = 91793.3
$ws.Range("J141").Value = 91793.3
$ws.Range("L141").Value = 91793.3
$ws.Range("N141").Value = -102153.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 63203756
$ws.Range("I4").Value = 84271020
$ws.Range("K4").Value = 252813060
$ws.Range("M4").Value = -252812948

$ws.Range("H37").Value = 113333.336
$ws.Range("J37").Value = 113333.336
$ws.Range("L37").Value = 340000.008
$ws.Range("N37").Value = -340224.008

$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").ClearContents()

$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").ClearContents()

$ws.Range("H113").Value = 1329.125
$ws.Range("I113").Value = 1873.75
$ws.Range("J113").Value = 1147.5834
$ws.Range("K113").Value = 5621.25
$ws.Range("L113").Value = 3442.7502
$ws.Range("M113").Value = -3451.25
$ws.Range("N113").Value = -7782.7502

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 14218.75
$ws.Range("I29").Value = 16750
$ws.Range("J29").Value = 13375
$ws.Range("K29").Value = 16750
$ws.Range("L29").Value = 13375
$ws.Range("M29").Value = -16460
$ws.Range("N29").Value = -13955

$ws.Range("H70").Value = 5750
$ws.Range("I70").Value = 8000
$ws.Range("K70").Value = 8000
$ws.Range("M70").Value = -7730

$ws.Range("H73").Value = 5750
$ws.Range("I73").Value = 8000
$ws.Range("K73").Value = 8000
$ws.Range("M73").Value = -7064

$ws.Range("H80").Value = 1941.1666
$ws.Range("I80").Value = 670.7143
$ws.Range("K80").Value = 670.7143
$ws.Range("M80").Value = 327.2857

$ws.Range("H83").Value = 1941.1666
$ws.Range("I83").Value = 670.7143
$ws.Range("K83").Value = 3353.5715
$ws.Range("M83").Value = 1638.4285

$ws.Range("H132").Value = 2583.375
$ws.Range("I132").Value = 2049.7273
$ws.Range("J132").Value = 3757.4
$ws.Range("K132").Value = 6149.1819
$ws.Range("L132").Value = 11272.2
$ws.Range("M132").Value = -3619.1819
$ws.Range("N132").Value = -16332.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 9999
$ws.Range("I16").Value = 9999
$ws.Range("K16").Value = 9999
$ws.Range("M16").Value = -9829

$ws.Range("H22").Value = 4999.3335
$ws.Range("I22").Value = 4999
$ws.Range("K22").Value = 4999
$ws.Range("M22").Value = -4704

$ws.Range("H27").Value = 4999.3335
$ws.Range("I27").Value = 4999
$ws.Range("K27").Value = 4999
$ws.Range("M27").Value = -4892

$ws.Range("H55").Value = 554.8333
$ws.Range("I55").Value = 475
$ws.Range("K55").Value = 475
$ws.Range("M55").Value = -302

$ws.Range("H68").Value = 2263.5715
$ws.Range("I68").Value = 2330
$ws.Range("K68").Value = 2330
$ws.Range("M68").Value = -1581

$ws.Range("H71").Value = 2263.5715
$ws.Range("I71").Value = 2330
$ws.Range("K71").Value = 11650
$ws.Range("M71").Value = -7906

$ws.Range("H121").Value = 34999.5
$ws.Range("J121").Value = 34999.5
$ws.Range("L121").Value = 34999.5
$ws.Range("N121").Value = -38493.5

$ws.Range("H132").Value = 3655.9211
$ws.Range("I132").Value = 3129.5334
$ws.Range("J132").Value = 5629.875
$ws.Range("K132").Value = 9388.600199999999
$ws.Range("L132").Value = 16889.625
$ws.Range("M132").Value = -6858.600199999999
$ws.Range("N132").Value = -21949.625

$ws.Range("H136").Value = 3038.6
$ws.Range("I136").Value = 3363.5
$ws.Range("J136").Value = 1739
$ws.Range("K136").Value = 10090.5
$ws.Range("L136").Value = 5217
$ws.Range("M136").Value = -7540.5
$ws.Range("N136").Value = -10317

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8221.9
$ws.Range("I62").Value = 6874
$ws.Range("J62").Value = 8558.875
$ws.Range("K62").Value = 6874
$ws.Range("L62").Value = 8558.875
$ws.Range("M62").Value = -6250
$ws.Range("N62").Value = -9806.875

$ws.Range("H65").Value = 8221.9
$ws.Range("I65").Value = 6874
$ws.Range("J65").Value = 8558.875
$ws.Range("K65").Value = 34370
$ws.Range("L65").Value = 42794.375
$ws.Range("M65").Value = -31250
$ws.Range("N65").Value = -49034.375

$ws.Range("H76").Value = 100000
$ws.Range("J76").Value = 100000
$ws.Range("L76").Value = 100000
$ws.Range("N76").Value = -100630

$ws.Range("H79").Value = 100000
$ws.Range("J79").Value = 100000
$ws.Range("L79").Value = 100000
$ws.Range("N79").Value = -102184

$ws.Range("H81").Value = 1973.4
$ws.Range("I81").Value = 1466.75
$ws.Range("J81").Value = 4000
$ws.Range("K81").Value = 2933.5
$ws.Range("L81").Value = 8000
$ws.Range("M81").Value = -1872.5
$ws.Range("N81").Value = -10122

$ws.Range("H84").Value = 1973.4
$ws.Range("I84").Value = 1466.75
$ws.Range("J84").Value = 4000
$ws.Range("K84").Value = 14667.5
$ws.Range("L84").Value = 40000
$ws.Range("M84").Value = -9363.5
$ws.Range("N84").Value = -50608
